$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20-21 block ("Unused land and land in transition")
$ws.Range("D20").Value = "Unused and transitioning land"
$ws.Range("F20").Value = "Rehabilitating land"
$ws.Range("F21").Value = "Unused land"

# Row 48-51 block ("Land in transition")
$ws.Range("D48").Value = "Vacant and transitioning land"
$ws.Range("F48").Value = "Rehabilitating land"
$ws.Range("F49").Value = "Abandoned land"
$ws.Range("F50").Value = "Degraded land"
$ws.Range("F51").Value = "No defined use"

# Row 93-95 block ("Vacant and transitioning land")
$ws.Range("F93").Value = "Brownfield development"
$ws.Range("F95").Value = "Idle and derelict land"
